# Update ExcelService and Wine model to ensure data integrity, add logging and new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row (row 1) : A..M
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Millesime"
$ws.Range("B1").Value = "Cuvee"
$ws.Range("C1").Value = "Domaine"
$ws.Range("D1").Value = "Appellation"
$ws.Range("E1").Value = "Region"
$ws.Range("F1").Value = "Country"
$ws.Range("G1").Value = "Pricetobuy"
$ws.Range("H1").Value = "Pricetosell"
$ws.Range("I1").Value = "Cost"
$ws.Range("J1").Value = "Quantity"
$ws.Range("K1").Value = "Updated"
$ws.Range("L1").Value = "ID"
$ws.Range("M1").Value = "Supplier"

# Header formatting: bold + centered (reuses existing style for A-H,J-M; creates new bold+0.00 style for I)
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1:H1").HorizontalAlignment = -4108
$ws.Range("J1:M1").Font.Bold = $true
$ws.Range("J1:M1").HorizontalAlignment = -4108

$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Data rows 2-6
# ---------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = 2018
$ws.Range("B2").Value = """Vieilles Vignes"""
$ws.Range("C2").Value = "Geanthet Pansiot"
$ws.Range("D2").Value = "Gevrey-Chambertin"
$ws.Range("E2").Value = "Bourgogne"
$ws.Range("F2").Value = "France"
$ws.Range("G2").Value = 138
$ws.Range("H2").Value = 650
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 45666
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "Cave Royale"

# Row 3
$ws.Range("A3").Value = 2020
$ws.Range("B3").Value = "Terroir Chapelle"
$ws.Range("C3").Value = "Patrick Piuze"
$ws.Range("D3").Value = "Chablis"
$ws.Range("E3").Value = "Bourgogne"
$ws.Range("F3").Value = "France"
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 90
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 45662
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = "Cavesa"

# Row 4
$ws.Range("A4").Value = 2010
$ws.Range("B4").Value = "Brut Vintage - Magnum"
$ws.Range("C4").Value = "Dom Pérignon"
$ws.Range("D4").Value = "Champagne"
$ws.Range("E4").Value = "Champagne"
$ws.Range("F4").Value = "France"
$ws.Range("G4").Value = 300
$ws.Range("H4").Value = 1300
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 45635
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = "Cave Royale"

# Row 5
$ws.Range("A5").Value = 2009
$ws.Range("B5").Value = "Brut Vintage Rosé"
$ws.Range("C5").Value = "Dom Pérignon"
$ws.Range("D5").Value = "Champagne"
$ws.Range("E5").Value = "Champagne"
$ws.Range("F5").Value = "France"
$ws.Range("G5").Value = 350
$ws.Range("H5").Value = 1450
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 45635
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = "Cavesa"

# Row 6
$ws.Range("A6").Value = 2008
$ws.Range("B6").Value = "Brut Vintage Rosé"
$ws.Range("C6").Value = "Dom Pérignon"
$ws.Range("D6").Value = "Champagne"
$ws.Range("E6").Value = "Champagne"
$ws.Range("F6").Value = "France"
$ws.Range("G6").Value = 350
$ws.Range("H6").Value = 1450
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 45662
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = "Cave Royale"

# ---------------------------------------------------------------
# Cost formulas (column I), numeric format 0.00
# ---------------------------------------------------------------
$ws.Range("I2").Formula = "=G2/(H2/1.081)"
$ws.Range("I3").Formula = "=G3/(H3/1.081)"
$ws.Range("I4").Formula = "=G4/(H4/1.081)"
$ws.Range("I5").Formula = "=G5/(H5/1.081)"
$ws.Range("I6").Formula = "=G6/(H6/1.081)"

$ws.Range("I2:I6").HorizontalAlignment = -4108
$ws.Range("I2:I6").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Alignment (center) for remaining plain cells -> reuses existing style 1
# ---------------------------------------------------------------
$ws.Range("A2:H6").HorizontalAlignment = -4108
$ws.Range("J2:J6").HorizontalAlignment = -4108
$ws.Range("L2:M6").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# Updated column (K) date format - applied to a single cell first then
# propagated via copy/paste-special (formats only) to avoid the engine
# creating a duplicate style per cell when a brand-new built-in date
# format is applied directly to a multi-cell range.
# ---------------------------------------------------------------
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").NumberFormat = "mm-dd-yy"
$ws.Range("K2").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Column widths (best effort; this engine quantizes ColumnWidth to
# 1/6-character increments, so exact fractional widths from real
# Excel cannot be reproduced bit-for-bit - we target the closest
# achievable value).
# ---------------------------------------------------------------
$ws.Columns("D:F").ColumnWidth = 19.022135416666668
$ws.Columns("G:H").ColumnWidth = 10.592447916666666
$ws.Columns("I:I").ColumnWidth = 9.307291666666666
$ws.Columns("J:L").ColumnWidth = 10.592447916666666
$ws.Columns("M:M").ColumnWidth = 12.022135416666666

# ---------------------------------------------------------------
# Sheet view / dimension cosmetics
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("E16").Select()
